$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 10) - same header labels as row 2
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = "__init__"
$ws.Range("C10").Value = "build_codebook"
$ws.Range("D10").Value = "from_string"
$ws.Range("E10").Value = "__encode_tree"
$ws.Range("F10").Value = "unzip_tree"
$ws.Range("G10").Value = "huffman_encode"
$ws.Range("H10").Value = "huffman_decode"

# Row 11 - weight
$ws.Range("A11").Value = "weight"
$ws.Range("B11").Value = "C"
$ws.Range("D11").Value = "O"

# Row 12 - data
$ws.Range("A12").Value = "data"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "O"
$ws.Range("D12").Value = "O"
$ws.Range("E12").Value = "O"
$ws.Range("F12").Value = "O"
$ws.Range("H12").Value = "O"

# Row 13 - left
$ws.Range("A13").Value = "left"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "O"
$ws.Range("E13").Value = "O"
$ws.Range("H13").Value = "O"

# Row 14 - right
$ws.Range("A14").Value = "right"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "O"
$ws.Range("E14").Value = "O"
$ws.Range("H14").Value = "O"

# Row 15 - codebook
$ws.Range("A15").Value = "codebook"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "T"
$ws.Range("G15").Value = "O"

# Copy formatting from the original table (rows 2-7) to the new table (rows 10-15)
$ws.Range("A2:H7").Copy()
$ws.Range("A10:H15").PasteSpecial(-4122)  # xlPasteFormats

# Set the selection to J7 as in the final state
$ws.Range("J7").Select()
